$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 7 (Adb) gains a hyperlink, like the other certificate rows above it ---
# Hyperlinks.Add() forces Excel's built-in "Hyperlink" cell style onto the
# target cell, which would normally fork off a brand-new cellXf. Restore the
# original (already-present) hyperlink-style formatting afterwards by copying
# it from A3, which already carries that exact style, so the workbook keeps
# reusing the existing style slot instead of creating a new one.
$ws.Hyperlinks.Add($ws.Range("A7"), "https://software-testing.ru/edu/3-online/272-adb") | Out-Null
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- New row 8: another certificate entry, "Chrome DevTools" ---
# Copy row 6's plain (non-hyperlinked) formatting down to row 8 first ...
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A8:C8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
# ... then fill in the course name.
$ws.Range("A8").Value = "Chrome DevTools"

# Column A now needs an explicit width to comfortably fit "Chrome DevTools".
$ws.Columns.Item(1).ColumnWidth = 14.75

Write-Host "done"
